# Finalizando Desenvolvimento do Banco de Dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Plan1" to "Médicos"
$ws.Name = "Médicos"

# Column widths (values chosen so the engine's internal width quantization
# lands as close as possible to the target stored widths 3.42578125 /
# 12.7109375 / 16.7109375 / 7.28515625)
$ws.Columns.Item(1).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(4).ColumnWidth = 6.5

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "ID_USUARIO"
$ws.Range("C1").Value = "Especialidade"
$ws.Range("D1").Value = "Crm"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 54356

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 53452

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = 65463

# --- Styling -----------------------------------------------------------
# Build each distinct final style exactly once on an out-of-the-way scratch
# cell, then stamp it onto the real destination cells with a single
# Copy/PasteSpecial(formats) call per style. Doing so avoids leaving behind
# unused "in-between" cell formats in the style table (which is what would
# happen if Font/Alignment properties were set directly, one after another,
# on the already-final multi-cell ranges).

# Style 1: default font, centered (-> A1:B4)
$s1 = $ws.Range("Z1")
$s1.HorizontalAlignment = -4108
$s1.VerticalAlignment = -4108

# Style 2: Arial 12, centered (-> C1:D1)
$s2 = $ws.Range("Z2")
$s2.Font.Name = "Arial"
$s2.Font.Size = 12
$s2.HorizontalAlignment = -4108
$s2.VerticalAlignment = -4108

# Style 3: Arial 11, centered (-> C2:D3, D4)
$s3 = $ws.Range("Z3")
$s3.Font.Name = "Arial"
$s3.Font.Size = 11
$s3.HorizontalAlignment = -4108
$s3.VerticalAlignment = -4108

# Style 4: Arial 11, color 4A4A4A, centered (-> C4)
$s4 = $ws.Range("Z4")
$s4.Font.Name = "Arial"
$s4.Font.Size = 11
$s4.Font.Color = 4868682
$s4.HorizontalAlignment = -4108
$s4.VerticalAlignment = -4108

$xlPasteFormats = -4122

$s1.Copy()
$ws.Range("A1:B4").PasteSpecial($xlPasteFormats)

$s2.Copy()
$ws.Range("C1:D1").PasteSpecial($xlPasteFormats)

$s3.Copy()
$ws.Range("C2:D3").PasteSpecial($xlPasteFormats)
$ws.Range("D4").PasteSpecial($xlPasteFormats)

$s4.Copy()
$ws.Range("C4").PasteSpecial($xlPasteFormats)

# Clean up the scratch cells used to build the styles
$ws.Range("Z1:Z4").Clear()

# Page margins (left/right were nudged slightly off the 0.7" default)
$ws.PageSetup.LeftMargin = 50.35
$ws.PageSetup.RightMargin = 50.35

# Selection matches the final sqref in the diff
$ws.Range("A1:D4").Select()
